$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2024-07-17 Wednesday" "2024-07-18 Thursday"
Replace-Text "891÷3=" "637÷5="
Replace-Text "227÷3=" "337÷9="
Replace-Text "337÷6=" "379÷9="
Replace-Text "909÷4=" "692÷9="
Replace-Text "590÷8=" "913÷4="
Replace-Text "916÷9=" "628÷8="
Replace-Text "221÷2=" "975÷9="
Replace-Text "258÷6=" "116÷7="
Replace-Text "975÷2=" "386÷6="
Replace-Text "587÷2=" "804÷8="
Replace-Text "944÷7=" "998÷5="
Replace-Text "175÷4=" "404÷8="
Replace-Text "705÷5=" "836÷2="
Replace-Text "750÷9=" "653÷4="
Replace-Text "779÷5=" "637÷9="
Replace-Text "171÷4=" "300÷6="
Replace-Text "839÷9=" "843÷8="
Replace-Text "156÷7=" "430÷4="
Replace-Text "193÷7=" "211÷3="
Replace-Text "730÷2=" "323÷5="
Replace-Text "447÷4=" "676÷2="
Replace-Text "985÷9=" "161÷4="
Replace-Text "426÷8=" "996÷2="
Replace-Text "808÷8=" "412÷4="
Replace-Text "925÷5=" "240÷2="

Write-Output "Done applying replacements"
